# se modif data para regresion en pre prod R31
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the account number in G3
$ws.Range("G3").Value = 30652299071

# Update the related counter in M3
$ws.Range("M3").Value = 306

# Leave row 3 selected (whole row, A3:XFD3), matching the saved selection state
$ws.Range("A3").EntireRow.Select()
